$wb = $excel.ActiveWorkbook

# Duplicate the original sheet to create the new "Cam_ContribExpenditureLobbyist1" sheet,
# placed immediately after the original sheet.
$ws1 = $wb.Worksheets.Item("Cam_ContribExpenditureLobbyistT")
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)

$ws2 = $wb.Worksheets.Item($ws1.Index + 1)
$ws2.Name = "Cam_ContribExpenditureLobbyist1"

# Repoint the defined name at the new sheet's data range.
$name = $wb.Names.Item("Cam_ContribExpenditureLobbyistType")
$name.RefersTo = "='Cam_ContribExpenditureLobbyist1'!`$A`$1:`$C`$10"
